$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(4, "Chanvre", 2, 0.45),
    @(17, "Fourrure", 38, 2.2000000000000002),
    @(18, "Cuir", 6, 0.321),
    @(19, "Argent", 54, 0.152),
    @(25, "Or", 102, 0.202)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}

$ws.Range("D7").Select()
